{"js": "// Remove the \"UMLS Quick Start Guide\" hyperlink and merge the surrounding\n// text into a single run, per the commit:\n//   \" and to explore the \" + [UMLS Quick Start Guide] + \", and other training materials.\"\n//   -> \" and to explore other training materials.\"\nconst body = context.document.body;\n\nconst oldText = \" and to explore the UMLS Quick Start Guide, and other training materials.\";\nconst newText = \" and to explore other training materials.\";\n\nconst results = body.search(oldText, { matchCase: true, ignorePunct: false, ignoreSpace: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  // Replace the whole matched range (run text + hyperlink + trailing run)\n  // with a single plain-text run, which removes the hyperlink entirely.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n} else {\n  throw new Error('Target text \"' + oldText + '\" not found in document body.');\n}\n", "ps1": "# Remove the \"UMLS Quick Start Guide\" hyperlink and merge the surrounding\n# text into a single run, per the commit:\n#   \" and to explore the \" + [UMLS Quick Start Guide] + \", and other training materials.\"\n#   -> \" and to explore other training materials.\"\n$d = $word.ActiveDocument\n\n# 1) Locate the \"UMLS Quick Start Guide\" hyperlink and strip its hyperlink-ness\n#    (Hyperlinks.Item(...).Delete() removes the field/relationship and leaves the\n#    display text behind as plain text in the run).\n$targetHyperlink = $null\nfor ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {\n    $candidate = $d.Hyperlinks.Item($i)\n    if ($candidate.TextToDisplay -eq \"UMLS Quick Start Guide\") {\n        $targetHyperlink = $candidate\n        break\n    }\n}\nif ($targetHyperlink -ne $null) {\n    $targetHyperlink.Delete()\n}\n\n# 2) Delete the now-plain-text \"UMLS Quick Start Guide\" wording.\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"UMLS Quick Start Guide\"\nif ($rng.Find.Execute()) {\n    $rng.Delete()\n}\n\n# 3) Collapse the leftover \"explore the , and other\" wording down to\n#    \"explore other\" so the sentence reads \"... and to explore other training\n#    materials.\"\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Text = \"explore the , and other\"\nif ($rng2.Find.Execute()) {\n    $rng2.Text = \"explore other\"\n}\n"}
